# Weekly update: a new price record (week of 2023-04-25) is inserted for
# "Terminal La Palmera de La Serena - Poroto granado" right after the
# existing 2023-03-27 row, pushing the later historical rows down by one.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 132 (shifts old rows 132-135 down to 133-136,
# inheriting formatting such as the date style from the row above).
$ws.Rows.Item(132).Insert()

# Populate the newly inserted row with the new weekly record.
$ws.Cells.Item(132, 1).Value = 8
$ws.Cells.Item(132, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(132, 3).Value = "Coquimbo"
$ws.Cells.Item(132, 4).Value = 45041
$ws.Cells.Item(132, 5).Value = 4
$ws.Cells.Item(132, 6).Value = 100112030
$ws.Cells.Item(132, 7).Value = "Poroto granado"
$ws.Cells.Item(132, 8).Value = "Sin especificar"
$ws.Cells.Item(132, 9).Value = "Primera"
$ws.Cells.Item(132, 10).Value = 500
$ws.Cells.Item(132, 11).Value = 35000
$ws.Cells.Item(132, 12).Value = 36000
$ws.Cells.Item(132, 13).Value = 35500
$ws.Cells.Item(132, 14).Value = "`$/malla 25 kilos"
$ws.Cells.Item(132, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(132, 16).Value = 1420
$ws.Cells.Item(132, 17).Value = 25
$ws.Cells.Item(132, 18).Value = "Hortaliza"
